$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Turn the three reference-list items ("D3 Block", "Stack-Overflow",
#    "D3-Gallery") into hyperlinks pointing at the sites they refer to.
#    "D3 Block" additionally ends up split into three runs ("D3 Bl" /
#    "o" / "ck") in the target document, so we nudge that one with an
#    extra (no-op) formatting toggle on the middle character to force
#    Word to split the run.
# ------------------------------------------------------------------

function Add-ListHyperlink($displayText, $url) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $displayText) {
            # Re-derive a plain offset-based Range (rather than reusing the
            # Paragraph's own Range object) before handing it to
            # Hyperlinks.Add - this runtime resolves the insertion point
            # incorrectly when the paragraph's native Range is passed in
            # directly.
            $s = $p.Range.Start
            $e = $p.Range.End - 1
            $rng = $d.Range($s, $e)
            $d.Hyperlinks.Add($rng, $url, $null, $null, $displayText) | Out-Null
            return $rng
        }
    }
    return $null
}

Add-ListHyperlink "Stack-Overflow" "https://stackoverflow.com" | Out-Null
Add-ListHyperlink "D3-Gallery" "https://www.d3-graph-gallery.com" | Out-Null
$d3Range = Add-ListHyperlink "D3 Block" "https://bl.ocks.org"

if ($d3Range -ne $null) {
    # Force "D3 Block" -> "D3 Bl" / "o" / "ck" run split, matching the
    # authored edit, by re-applying (identical) direct formatting to the
    # single "o" character in the middle of the word.
    $s = $d3Range.Start
    $mid = $d.Range($s + 5, $s + 6)
    $mid.Font.Bold = 1
    $mid.Font.Bold = 0
}

# ------------------------------------------------------------------
# 2) Append the bare URL after each "... at :" reference paragraph.
# ------------------------------------------------------------------

function Append-UrlToParagraph($containsText, $url) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Contains($containsText)) {
            $endPos = $p.Range.End - 1
            $ins = $d.Range($endPos, $endPos)
            $ins.InsertAfter(" ")
            $ins2 = $d.Range($endPos + 1, $endPos + 1)
            $ins2.InsertAfter($url)
            return
        }
    }
}

Append-UrlToParagraph "Block.org" "https://bl.ocks.org"
Append-UrlToParagraph "[3] Stack-Overflow" "https://stackoverflow.com"
Append-UrlToParagraph "[4] D3-Gallery" "https://www.d3-graph-gallery.com"

Write-Host "References added"
